$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet/tab
$ws.Name = "iteration1"

# 2. Highlight the header rows (2-3) with the new fill + bold font
#    (Accent3, Lighter 60% -> #D7E4BD; Bold, size 13)
$header = $ws.Range("A2:K3")
$header.Font.Bold = $true
$header.Font.Size = 13
$header.Font.Name = "Calibri"
$header.Interior.Color = 12444887

# 3. Fill in the previously empty "Man Hours" (E) and extra (F) metrics
#    for the existing feature rows
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 0

$ws.Range("E5").Value = 14
$ws.Range("F5").Value = 0

$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 0

$ws.Range("E7").Value = 12
$ws.Range("F7").Value = 0

$ws.Range("E8").Value = 18
$ws.Range("F8").Value = 0

# 4. Row 6 ("View of Complete features provided by the system") now has
#    its metrics filled in as well
$ws.Range("C9").Formula = "= 149 + 239"
$ws.Range("D9").Formula = "= 0 + 8"
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 0
$ws.Range("G9").Value2 = $ws.Range("G8").Value2

# 5. Update the saved selection to match the author's cursor position
$ws.Range("F8").Select()
